$d = $word.ActiveDocument

# --- Paragraph: "cool" / " story bro " (split across two runs, separated by
#     proofErr gramStart/gramEnd markers) -> a single run "cool story bro "
#     with the proofErr markers removed.
$p1 = $d.Paragraphs(4)
$r1 = $p1.Range
[void]$r1.MoveEnd(1, -1)           # exclude the paragraph mark
$r1.InsertBefore("Z")              # pushes any leading proofErr marker to sit
                                    # inside the (now live/expanded) range
$r1.Text = "cool story bro "       # collapses every run + enclosed proofErr
                                    # marker in the range into one clean run

# --- Paragraph: "(" / "shrug" / ")" (three runs, proofErr gramStart/gramEnd
#     between them) -> a single run "(shrug)" with no proofErr markers.
$p2 = $d.Paragraphs(7)
$r2 = $p2.Range
[void]$r2.MoveEnd(1, -1)
$r2.InsertBefore("Z")
$r2.Text = "(shrug)"

# --- Paragraph containing the _GoBack bookmark -> plain empty paragraph.
$d.Bookmarks("_GoBack").Delete()
